# "mein name hat gefehlt" - add the missing name "Jonathan" to the list of
# names on the title slide (slide 1, subtitle placeholder "Untertitel 2").
#
# The subtitle text frame currently reads (with tab-separated columns):
#   Ausgearbeitet von
#   Pablo                    Anna
#   „“                       Paula      <- placeholder quotes where a name is missing
#   Fabian                   ... und einigen Mentoren
#
# We replace the empty-name placeholder quotes „“ with "Jonathan" and trim
# two tabs from the run so "Paula" stays reasonably aligned with the other
# names, giving:
#   Jonathan                 Paula

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)          # "Untertitel 2" subtitle placeholder
$tr = $sh.TextFrame.TextRange

# Locate "Paula" - this is a stable anchor immediately after the placeholder
# quotes + tabs, and lets us compute the placeholder's position without
# relying on literal smart-quote characters.
$f = $tr.Find("Paula")
$paulaStart = $f.Start

# The run reads: „ “ (2 chars) + 6 tabs + "Paula", so the placeholder quotes
# sit 8 characters before "Paula" begins.
$quotes = $tr.Characters($paulaStart - 8, 2)
$quotes.Text = "Jonathan"

# Re-find "Paula" (its start shifted after the text above grew) and drop two
# of the six tabs right before it, so the name lines up with the other rows.
$f2 = $tr.Find("Paula")
$paulaStart2 = $f2.Start
$extraTabs = $tr.Characters($paulaStart2 - 2, 2)
$extraTabs.Text = ""
